$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text values (e.g. "67.766.14", "1.00").
# Assigning a numeric-looking string straight to .Value lets Excel coerce it
# to a real number, which would change the cell's stored type. Prefixing the
# literal with an apostrophe via .Formula keeps it text, exactly like typing
# '67.766.14 into the cell in the Excel UI.
$ws.Range("D2").Formula  = "'67.649.91"
$ws.Range("E2").Value    = "  +0.04%  "

$ws.Range("D3").Formula  = "'2.614.96"
$ws.Range("E3").Value    = "  -1.04%  "

$ws.Range("E4").Value    = "  +0.02%  "

$ws.Range("D5").Formula  = "'594.33"
$ws.Range("E5").Value    = "  -1.89%  "

$ws.Range("D6").Formula  = "'152.12"
$ws.Range("E6").Value    = "  -1.65%  "

$ws.Range("E7").Value    = "  +0.04%  "

$ws.Range("D8").Formula  = "'0.544"
$ws.Range("E8").Value    = "  -0.87%  "

$ws.Range("D9").Formula  = "'2.613.71"
$ws.Range("E9").Value    = "  -1.05%  "

$ws.Range("E10").Value   = "  +5.77%  "

$ws.Range("E11").Value   = "  -0.73%  "

$ws.Range("D12").Formula = "'5.18"
$ws.Range("E12").Value   = "  -0.76%  "

$ws.Range("E13").Value   = "  -1.79%  "

$ws.Range("D14").Formula = "'27.45"
$ws.Range("E14").Value   = "  -2.24%  "

$ws.Range("D15").Formula = "'0.0000187"
$ws.Range("E15").Value   = "  +1.34%  "

$ws.Range("D16").Formula = "'3.097.24"
$ws.Range("E16").Value   = "  -0.78%  "

$ws.Range("D17").Formula = "'67.664.21"
$ws.Range("E17").Value   = "  +0.12%  "

$ws.Range("D18").Formula = "'2.610.08"
$ws.Range("E18").Value   = "  -1.09%  "

$ws.Range("D19").Formula = "'371.30"
$ws.Range("E19").Value   = "  +1.53%  "

$ws.Range("D20").Formula = "'11.19"
$ws.Range("E20").Value   = "  -1.36%  "

$ws.Range("D21").Formula = "'7.41"
$ws.Range("E21").Value   = "  -3.03%  "

$ws.Range("D22").Formula = "'4.22"
$ws.Range("E22").Value   = "  -2.20%  "

$ws.Range("D23").Formula = "'4.80"
$ws.Range("E23").Value   = "  -3.71%  "

$ws.Range("D24").Formula = "'2.04"
$ws.Range("E24").Value   = "  -5.28%  "

$ws.Range("D25").Formula = "'72.57"
$ws.Range("E25").Value   = "  +9.74%  "

$ws.Range("D26").Formula = "'0.999"
$ws.Range("E26").Value   = "  -0.16%  "

$ws.Range("D27").Formula = "'9.82"
$ws.Range("E27").Value   = "  -2.17%  "

$ws.Range("D28").Formula = "'594.71"
$ws.Range("E28").Value   = "  +1.93%  "

$ws.Range("D29").Formula = "'2.750.02"
$ws.Range("E29").Value   = "  -0.32%  "

$ws.Range("E30").Value   = "  -1.99%  "

$ws.Range("D31").Formula = "'1.00"
$ws.Range("E31").Value   = "  +0.10%  "

$ws.Range("D32").Formula = "'7.78"
$ws.Range("E32").Value   = "  -2.00%  "

$ws.Range("D33").Formula = "'1.37"
$ws.Range("E33").Value   = "  -3.76%  "

$ws.Range("D34").Formula = "'1.84"
$ws.Range("E34").Value   = "  -0.96%  "

$ws.Range("D35").Formula = "'1.00"
$ws.Range("E35").Value   = "  +0.03%  "

$ws.Range("E36").Value   = "  -3.85%  "

$ws.Range("E37").Value   = "  -1.87%  "

$ws.Range("D38").Formula = "'158.06"
$ws.Range("E38").Value   = "  +0.01%  "

$ws.Range("D39").Formula = "'19.09"
$ws.Range("E39").Value   = "  -2.17%  "

$ws.Range("E40").Value   = "  +2.42%  "

$ws.Range("D41").Formula = "'0.367"
$ws.Range("E41").Value   = "  -1.54%  "

$ws.Range("D42").Formula = "'5.25"
$ws.Range("E42").Value   = "  -1.69%  "

$ws.Range("D43").Formula = "'2.67"
$ws.Range("E43").Value   = "  +1.38%  "

$ws.Range("E44").Value   = "  +4.47%  "

$ws.Range("E45").Value   = "  +0.02%  "

$ws.Range("D46").Formula = "'40.38"
$ws.Range("E46").Value   = "  -1.92%  "

$ws.Range("D47").Formula = "'155.83"
$ws.Range("E47").Value   = "  -0.60%  "

$ws.Range("D48").Formula = "'0.0₆0295"
$ws.Range("E48").Value   = "  +1.69%  "

$ws.Range("D49").Formula = "'3.67"
$ws.Range("E49").Value   = "  -1.96%  "

$ws.Range("D50").Formula = "'1.68"
$ws.Range("E50").Value   = "  -3.23%  "

$ws.Range("D51").Formula = "'0.0777"
$ws.Range("E51").Value   = "  -1.58%  "
